$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistic (column C) and p-value (column D) figures
# for rows 2-11, per the corrected computation.

$ws.Range("C2").Value = 0.7507613260376965
$ws.Range("D2").Value = 0.4607516456569449

$ws.Range("C3").Value = -2.336156315807316
$ws.Range("D3").Value = 0.02899750939657242

$ws.Range("C4").Value = -1.297421994231022
$ws.Range("D4").Value = 0.2079267239804898

$ws.Range("C5").Value = -3.35773819112655
$ws.Range("D5").Value = 0.00284382980167952

$ws.Range("C6").Value = -2.428657215326825
$ws.Range("D6").Value = 0.02378098985670718

$ws.Range("C7").Value = -1.29897213335943
$ws.Range("D7").Value = 0.2074034518480614

$ws.Range("C8").Value = -3.221075724846067
$ws.Range("D8").Value = 0.003931772905708897

$ws.Range("C9").Value = 0.4969166119581448
$ws.Range("D9").Value = 0.6241772515538524

$ws.Range("C10").Value = -2.72614844384788
$ws.Range("D10").Value = 0.01233208724704737

$ws.Range("C11").Value = -2.153339967598714
$ws.Range("D11").Value = 0.04251475039955421
